# "Create and delete data in postgres and orion"
#
# The "category" sample sheet moves from a single "#テナント名" (tenant
# name) column to two columns - "テナント名" and a new "サービスパス"
# (service path) column - and the "政策区域" (government-plan) category
# row is removed entirely while its tenant identifiers are replaced with
# the new Saitama/Sakado tenant id and per-category service paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the old A1:F4 block so no stale cells are
# left behind once the new layout (7 cols x 3 rows) is written.
$ws.Range("A1:F4").ClearContents()

# --- Header row ---------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "#カテゴリID"
$ws.Cells.Item(1, 2).Value = "カテゴリ名"
$ws.Cells.Item(1, 3).Value = "テナント名"
$ws.Cells.Item(1, 4).Value = "サービスパス"
$ws.Cells.Item(1, 5).Value = "色"
$ws.Cells.Item(1, 6).Value = "表示順"
$ws.Cells.Item(1, 7).Value = "有効"

# --- Row 2: 公共施設 (public facility) ----------------------------------
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "公共施設"
$ws.Cells.Item(2, 3).Value = "112399_sakado_city"
$ws.Cells.Item(2, 4).Value = "/public_facility"
$ws.Cells.Item(2, 5).Value = "#00008b"
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = "〇"

# --- Row 3: 水域情報 (flood control) -------------------------------------
# Category id 3 keeps its original id (the 政策区域/government-plan row,
# id 2, is dropped), so it lands on sheet row 3 once the old row 3 is
# removed below.
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "水域情報"
$ws.Cells.Item(3, 3).Value = "112399_sakado_city"
$ws.Cells.Item(3, 4).Value = "/flood_control"
$ws.Cells.Item(3, 5).Value = "#65ace4"
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = "〇"

# The 政策区域 (government plan) category is deleted outright - remove the
# now-stale leftover row (old row 4) and shift everything below it up.
$ws.Rows.Item(4).Delete()

# New サービスパス column D inherits column C's width, same as the other
# widened columns.
$ws.Columns.Item(4).ColumnWidth = 30.43

# Match the saved selection left behind in the source workbook.
$ws.Range("A5").Select() | Out-Null
